$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# --- Rename existing "speed_up" column (G) to "speed_up_init" ---
$ws.Range("G1").Value = "speed_up_init"

# --- Add two new calculated columns to the table: speed_up_add (H), speed_up_tot (I) ---
$colAdd = $tbl.ListColumns.Add()
$colAdd.Range.Item(1).Value = "speed_up_add"

$colTot = $tbl.ListColumns.Add()
$colTot.Range.Item(1).Value = "speed_up_tot"

# --- Fill in formulas row by row (avoids shared-formula compression so the
#     saved XML keeps one <f> per cell, matching the source workbook's style) ---
for ($r = 2; $r -le 13; $r++) {
  $ws.Range("G$r").Formula = "=1278/Tabella1[[#This Row],[init_min_time]]"
  $ws.Range("H$r").Formula = "=294/Tabella1[[#This Row],[add_min_time]]"
  $ws.Range("I$r").Formula = "=1572/Tabella1[[#This Row],[total_time]]"
}

# --- Number formats for the new speed-up columns (integer display) ---
$ws.Range("G2:I13").NumberFormat = "0"
$ws.Range("I10").HorizontalAlignment = -4152

# --- Column widths (best effort; engine snaps to coarser character grid) ---
$ws.Columns.Item(6).ColumnWidth = $ws.Columns.Item(5).ColumnWidth
$ws.Columns.Item(7).ColumnWidth = 13.94
$ws.Columns.Item(8).ColumnWidth = 16.28
$ws.Columns.Item(9).ColumnWidth = 15.61

# --- Selection moves to the newly added columns ---
$null = $ws.Range("G1:I13").Select()

$wb.Save()
